# last test and twist test
# Adds an 11th participant/column (L) of results to the "Questionnarie"
# and "Time_experiment1.1" sheets, and updates each sheet's selection.

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------
# Sheet "Questionnarie" (sheet1) - new column L (12), rows 1-30
# ---------------------------------------------------------------------
$ws1 = $wb.Worksheets.Item("Questionnarie")

$ws1.Cells.Item(1, 12).Value = 11

$ws1.Cells.Item(2, 12).Value = "F"
$ws1.Cells.Item(3, 12).Value = 37
$ws1.Cells.Item(4, 12).Value = "N"
$ws1.Cells.Item(5, 12).Value = "N"
$ws1.Cells.Item(6, 12).Value = 4
$ws1.Cells.Item(7, 12).Value = "Y"
$ws1.Cells.Item(8, 12).Value = "N"
$ws1.Cells.Item(9, 12).Value = "Y"
$ws1.Cells.Item(10, 12).Value = 3
$ws1.Cells.Item(11, 12).Value = 3
$ws1.Cells.Item(12, 12).Value = 1
$ws1.Cells.Item(13, 12).Value = 5
$ws1.Cells.Item(14, 12).Value = 1
$ws1.Cells.Item(15, 12).Value = 2
$ws1.Cells.Item(16, 12).Value = 2
$ws1.Cells.Item(17, 12).Value = 1
$ws1.Cells.Item(18, 12).Value = 2
$ws1.Cells.Item(19, 12).Value = 5
$ws1.Cells.Item(20, 12).Value = 3
$ws1.Cells.Item(21, 12).Value = 1
$ws1.Cells.Item(22, 12).Value = 3
$ws1.Cells.Item(23, 12).Value = 2
$ws1.Cells.Item(24, 12).Value = 3
$ws1.Cells.Item(25, 12).Value = 3
$ws1.Cells.Item(26, 12).Value = 5
$ws1.Cells.Item(27, 12).Value = 5
$ws1.Cells.Item(28, 12).Value = 5
$ws1.Cells.Item(29, 12).Value = 1
$ws1.Cells.Item(30, 12).Value = 3

# Rows 2-30 share the centered style already used by columns B:K.
$ws1.Range("L2:L30").HorizontalAlignment = -4108

# Restore the sheet's recorded selection.
$ws1.Range("O21").Select()

# ---------------------------------------------------------------------
# Sheet "Time_experiment1.1" (sheet2) - new column L (12), rows 1-11
# ---------------------------------------------------------------------
$ws2 = $wb.Worksheets.Item("Time_experiment1.1")

$ws2.Cells.Item(1, 12).Value = 11
$ws2.Cells.Item(2, 12).Value = 11
$ws2.Cells.Item(3, 12).Value = 11
$ws2.Cells.Item(4, 12).Value = 11
$ws2.Cells.Item(5, 12).Value = 11
$ws2.Cells.Item(6, 12).Value = 12
$ws2.Cells.Item(7, 12).Value = 13
$ws2.Cells.Item(8, 12).Value = 13
$ws2.Cells.Item(9, 12).Value = 11
$ws2.Cells.Item(10, 12).Value = 11
$ws2.Cells.Item(11, 12).Value = 9

# Restore the sheet's recorded selection. Done last so this sheet
# (which was already the active tab) remains the active tab on save.
$ws2.Range("N14").Select()
